$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F ("Team Name") before the existing "Name" column,
# shifting F:J (Name, Order, Split, Leg, Cumulative) to G:K.
$ws.Columns("F:F").Insert()

$ws.Range("F1").Value = "Team Name"

$teams = @(
    "Tennessee","Tennessee","Tennessee","Tennessee","ASU","ASU","ASU","ASU","Florida","Florida",
    "Florida","Florida","California","California","California","California","Indiana","Indiana","Indiana","Indiana",
    "NC State","NC State","NC State","NC State","Texas","Texas","Texas","Texas","VT","VT",
    "VT","VT","Alabama","Alabama","Alabama","Alabama","Stanford","Stanford","Stanford","Stanford",
    "Georgia","Georgia","Georgia","Georgia","Arizona","Arizona","Arizona","Arizona","OSU","OSU",
    "OSU","OSU","TAMU","TAMU","TAMU","TAMU","Louisville","Louisville","Louisville","Louisville",
    "Virginia","Virginia","Virginia","Virginia","FSU","FSU","FSU","FSU","Michigan","Michigan",
    "Michigan","Michigan","Wisconsin","Wisconsin","Wisconsin","Wisconsin","Auburn","Auburn","Auburn","Auburn",
    "UNC","UNC","UNC","UNC","LSU","LSU","LSU","LSU","Yale","Yale",
    "Yale","Yale","Northwestern","Northwestern","Northwestern","Northwestern","USC","USC","USC","USC",
    "SMU","SMU","SMU","SMU","Princeton","Princeton","Princeton","Princeton","GT","GT",
    "GT","GT","Missouri","Missouri","Missouri","Missouri","Harvard","Harvard","Harvard","Harvard"
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $teams[$i]
}

Write-Host "Done. Dimension now:" $ws.UsedRange.Address()
